$d = $word.ActiveDocument
$d.Content.Find.Execute("loaded date to get", $true, $false, $false, $false, $false,
                         $true, 1, $false, "loaded data to get", 2)
